$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (A:T) reflecting refreshed TPM values.
# Columns: A=Sending cluster, B=Ligand symbol, C=Receptor symbol, D=Target cluster,
# E..T = numeric metric columns.
$rows = @(
    @("FAPs","Efna2","Epha5","ECs",3,1,2.015377,6.046131,0.7554960962715589,0.7554960962715588,1,0.3333333333333333,0.008504666666666666,0.025514,0.3779124020559003,0.3779124020559003,0.01714010959266666,0.154260986334,0.2855113444858405,0.2855113444858405),
    @("FAPs","Efna2","Epha5","MuSCs",3,1,2.015377,6.046131,0.7554960962715589,0.7554960962715588,1,0.3333333333333333,0.01399966666666667,0.041999,0.6220875979440996,0.6220875979440996,0.02821460620766667,0.253931455869,0.4699847517857183,0.4699847517857182),
    @("MuSCs","Efna2","Epha5","ECs",3,1,0.6522436666666667,1.956731,0.2445039037284412,0.2445039037284411,1,0.3333333333333333,0.008504666666666666,0.025514,0.3779124020559003,0.3779124020559003,0.005547114970444444,0.049924034734,0.09240105757005979,0.09240105757005979),
    @("MuSCs","Efna2","Epha5","MuSCs",3,1,0.6522436666666667,1.956731,0.2445039037284412,0.2445039037284411,1,0.3333333333333333,0.01399966666666667,0.041999,0.6220875979440996,0.6220875979440996,0.009131193918777779,0.082180745269,0.1521028461583813,0.1521028461583813)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowData[$c]
    }
}

$wb.Save()
